# "contingencies with rene fine"
#
# Insert two brand-new contingency rows ("line7", "line8") right after the
# existing "line6" row. This pushes the eight "extr*" rows down by two rows
# (row 8->10 .. row 15->17) and extends the sheet from A1:E15 to A1:E17.
# Besides the insertion, a handful of the (now shifted) extr* rows get their
# from_bus/to_bus/in_service values refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing extr1..extr8 rows (old rows 8-15) down to rows 10-17,
# preserving their values/format, and opening up rows 8-9 for the new lines.
$ws.Rows("8:9").Insert()

# --- New row 8: line7 --------------------------------------------------
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# --- New row 9: line8 --------------------------------------------------
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# Give the two new index cells (column A) the same look as the rest of the
# index column (bold, centered, bordered).
$ws.Range("A10").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Renumber the shifted index column (A) for rows 10-17 --------------
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(17, 1).Value = 15

# --- Refresh the few C/D/E values that changed on the shifted rows -----
# row 10 = extr1
$ws.Cells.Item(10, 5).Value = $true
# row 11 = extr2 (unchanged)
# row 12 = extr3 (unchanged)
# row 13 = extr4
$ws.Cells.Item(13, 5).Value = $true
# row 14 = extr5
$ws.Cells.Item(14, 5).Value = $true
# row 15 = extr6
$ws.Cells.Item(15, 5).Value = $true
# row 16 = extr7 (unchanged)
# row 17 = extr8
$ws.Cells.Item(17, 5).Value = $false
